$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 60,20
$data[0,0] = 2
$data[0,1] = 'Comercializadora del Agro de Limarí'
$data[0,2] = 'Coquimbo'
$data[0,3] = 44952
$data[0,4] = 4
$data[0,5] = 'Fruta'
$data[0,6] = 100103
$data[0,7] = 'Frutos de hueso (carozo)'
$data[0,8] = 100103002
$data[0,9] = 'Ciruela'
$data[0,10] = 'Black Amber'
$data[0,11] = 'Primera'
$data[0,12] = 10
$data[0,13] = 300000
$data[0,14] = 310000
$data[0,15] = 305000
$data[0,16] = '$/bins (450 kilos)'
$data[0,17] = 'Región de O''Higgins'
$data[0,18] = 678
$data[0,19] = 450
$data[1,0] = 2
$data[1,1] = 'Comercializadora del Agro de Limarí'
$data[1,2] = 'Coquimbo'
$data[1,3] = 44952
$data[1,4] = 4
$data[1,5] = 'Fruta'
$data[1,6] = 100103
$data[1,7] = 'Frutos de hueso (carozo)'
$data[1,8] = 100103002
$data[1,9] = 'Ciruela'
$data[1,10] = 'Black Amber'
$data[1,11] = 'Segunda'
$data[1,12] = 10
$data[1,13] = 230000
$data[1,14] = 240000
$data[1,15] = 235000
$data[1,16] = '$/bins (450 kilos)'
$data[1,17] = 'Región de O''Higgins'
$data[1,18] = 522
$data[1,19] = 450
$data[2,0] = 2
$data[2,1] = 'Comercializadora del Agro de Limarí'
$data[2,2] = 'Coquimbo'
$data[2,3] = 44615
$data[2,4] = 4
$data[2,5] = 'Fruta'
$data[2,6] = 100103
$data[2,7] = 'Frutos de hueso (carozo)'
$data[2,8] = 100103002
$data[2,9] = 'Ciruela'
$data[2,10] = 'Angeleno'
$data[2,11] = 'Primera'
$data[2,12] = 16
$data[2,13] = 200000
$data[2,14] = 210000
$data[2,15] = 205000
$data[2,16] = '$/bins (450 kilos)'
$data[2,17] = 'Región Metropolitana'
$data[2,18] = 456
$data[2,19] = 450
$data[3,0] = 2
$data[3,1] = 'Comercializadora del Agro de Limarí'
$data[3,2] = 'Coquimbo'
$data[3,3] = 44615
$data[3,4] = 4
$data[3,5] = 'Fruta'
$data[3,6] = 100103
$data[3,7] = 'Frutos de hueso (carozo)'
$data[3,8] = 100103002
$data[3,9] = 'Ciruela'
$data[3,10] = 'Angeleno'
$data[3,11] = 'Segunda'
$data[3,12] = 20
$data[3,13] = 160000
$data[3,14] = 170000
$data[3,15] = 165000
$data[3,16] = '$/bins (450 kilos)'
$data[3,17] = 'Región Metropolitana'
$data[3,18] = 367
$data[3,19] = 450
$data[4,0] = 2
$data[4,1] = 'Comercializadora del Agro de Limarí'
$data[4,2] = 'Coquimbo'
$data[4,3] = 44594
$data[4,4] = 4
$data[4,5] = 'Fruta'
$data[4,6] = 100103
$data[4,7] = 'Frutos de hueso (carozo)'
$data[4,8] = 100103002
$data[4,9] = 'Ciruela'
$data[4,10] = 'Black Amber'
$data[4,11] = 'Especial'
$data[4,12] = 240
$data[4,13] = 15500
$data[4,14] = 16000
$data[4,15] = 15750
$data[4,16] = '$/caja 15 kilos granel'
$data[4,17] = 'Región de O''Higgins'
$data[4,18] = 1050
$data[4,19] = 15
$data[5,0] = 2
$data[5,1] = 'Comercializadora del Agro de Limarí'
$data[5,2] = 'Coquimbo'
$data[5,3] = 44594
$data[5,4] = 4
$data[5,5] = 'Fruta'
$data[5,6] = 100103
$data[5,7] = 'Frutos de hueso (carozo)'
$data[5,8] = 100103002
$data[5,9] = 'Ciruela'
$data[5,10] = 'Black Amber'
$data[5,11] = 'Primera'
$data[5,12] = 300
$data[5,13] = 13500
$data[5,14] = 14000
$data[5,15] = 13750
$data[5,16] = '$/caja 15 kilos granel'
$data[5,17] = 'Región de O''Higgins'
$data[5,18] = 917
$data[5,19] = 15
$data[6,0] = 2
$data[6,1] = 'Comercializadora del Agro de Limarí'
$data[6,2] = 'Coquimbo'
$data[6,3] = 44720
$data[6,4] = 4
$data[6,5] = 'Fruta'
$data[6,6] = 100103
$data[6,7] = 'Frutos de hueso (carozo)'
$data[6,8] = 100103002
$data[6,9] = 'Ciruela'
$data[6,10] = 'Angeleno'
$data[6,11] = 'Especial'
$data[6,12] = 16
$data[6,13] = 220000
$data[6,14] = 230000
$data[6,15] = 225000
$data[6,16] = '$/bins (450 kilos)'
$data[6,17] = 'Región de O''Higgins'
$data[6,18] = 500
$data[6,19] = 450
$data[7,0] = 2
$data[7,1] = 'Comercializadora del Agro de Limarí'
$data[7,2] = 'Coquimbo'
$data[7,3] = 44720
$data[7,4] = 4
$data[7,5] = 'Fruta'
$data[7,6] = 100103
$data[7,7] = 'Frutos de hueso (carozo)'
$data[7,8] = 100103002
$data[7,9] = 'Ciruela'
$data[7,10] = 'Angeleno'
$data[7,11] = 'Primera'
$data[7,12] = 20
$data[7,13] = 190000
$data[7,14] = 200000
$data[7,15] = 195000
$data[7,16] = '$/bins (450 kilos)'
$data[7,17] = 'Región de O''Higgins'
$data[7,18] = 433
$data[7,19] = 450
$data[8,0] = 2
$data[8,1] = 'Comercializadora del Agro de Limarí'
$data[8,2] = 'Coquimbo'
$data[8,3] = 44720
$data[8,4] = 4
$data[8,5] = 'Fruta'
$data[8,6] = 100103
$data[8,7] = 'Frutos de hueso (carozo)'
$data[8,8] = 100103002
$data[8,9] = 'Ciruela'
$data[8,10] = 'Angeleno'
$data[8,11] = 'Segunda'
$data[8,12] = 20
$data[8,13] = 150000
$data[8,14] = 160000
$data[8,15] = 155000
$data[8,16] = '$/bins (450 kilos)'
$data[8,17] = 'Región de O''Higgins'
$data[8,18] = 344
$data[8,19] = 450
$data[9,0] = 2
$data[9,1] = 'Comercializadora del Agro de Limarí'
$data[9,2] = 'Coquimbo'
$data[9,3] = 44679
$data[9,4] = 4
$data[9,5] = 'Fruta'
$data[9,6] = 100103
$data[9,7] = 'Frutos de hueso (carozo)'
$data[9,8] = 100103002
$data[9,9] = 'Ciruela'
$data[9,10] = 'Angeleno'
$data[9,11] = 'Especial'
$data[9,12] = 10
$data[9,13] = 220000
$data[9,14] = 230000
$data[9,15] = 225000
$data[9,16] = '$/bins (450 kilos)'
$data[9,17] = 'Región de O''Higgins'
$data[9,18] = 500
$data[9,19] = 450
$data[10,0] = 2
$data[10,1] = 'Comercializadora del Agro de Limarí'
$data[10,2] = 'Coquimbo'
$data[10,3] = 44679
$data[10,4] = 4
$data[10,5] = 'Fruta'
$data[10,6] = 100103
$data[10,7] = 'Frutos de hueso (carozo)'
$data[10,8] = 100103002
$data[10,9] = 'Ciruela'
$data[10,10] = 'Angeleno'
$data[10,11] = 'Primera'
$data[10,12] = 20
$data[10,13] = 180000
$data[10,14] = 190000
$data[10,15] = 185000
$data[10,16] = '$/bins (450 kilos)'
$data[10,17] = 'Región de O''Higgins'
$data[10,18] = 411
$data[10,19] = 450
$data[11,0] = 2
$data[11,1] = 'Comercializadora del Agro de Limarí'
$data[11,2] = 'Coquimbo'
$data[11,3] = 44679
$data[11,4] = 4
$data[11,5] = 'Fruta'
$data[11,6] = 100103
$data[11,7] = 'Frutos de hueso (carozo)'
$data[11,8] = 100103002
$data[11,9] = 'Ciruela'
$data[11,10] = 'Angeleno'
$data[11,11] = 'Segunda'
$data[11,12] = 20
$data[11,13] = 150000
$data[11,14] = 160000
$data[11,15] = 155000
$data[11,16] = '$/bins (450 kilos)'
$data[11,17] = 'Región de O''Higgins'
$data[11,18] = 344
$data[11,19] = 450
$data[12,0] = 2
$data[12,1] = 'Comercializadora del Agro de Limarí'
$data[12,2] = 'Coquimbo'
$data[12,3] = 44665
$data[12,4] = 4
$data[12,5] = 'Fruta'
$data[12,6] = 100103
$data[12,7] = 'Frutos de hueso (carozo)'
$data[12,8] = 100103002
$data[12,9] = 'Ciruela'
$data[12,10] = 'Angeleno'
$data[12,11] = 'Especial'
$data[12,12] = 16
$data[12,13] = 230000
$data[12,14] = 240000
$data[12,15] = 235000
$data[12,16] = '$/bins (450 kilos)'
$data[12,17] = 'Región Metropolitana'
$data[12,18] = 522
$data[12,19] = 450
$data[13,0] = 2
$data[13,1] = 'Comercializadora del Agro de Limarí'
$data[13,2] = 'Coquimbo'
$data[13,3] = 44665
$data[13,4] = 4
$data[13,5] = 'Fruta'
$data[13,6] = 100103
$data[13,7] = 'Frutos de hueso (carozo)'
$data[13,8] = 100103002
$data[13,9] = 'Ciruela'
$data[13,10] = 'Angeleno'
$data[13,11] = 'Primera'
$data[13,12] = 16
$data[13,13] = 190000
$data[13,14] = 200000
$data[13,15] = 195000
$data[13,16] = '$/bins (450 kilos)'
$data[13,17] = 'Región Metropolitana'
$data[13,18] = 433
$data[13,19] = 450
$data[14,0] = 2
$data[14,1] = 'Comercializadora del Agro de Limarí'
$data[14,2] = 'Coquimbo'
$data[14,3] = 44665
$data[14,4] = 4
$data[14,5] = 'Fruta'
$data[14,6] = 100103
$data[14,7] = 'Frutos de hueso (carozo)'
$data[14,8] = 100103002
$data[14,9] = 'Ciruela'
$data[14,10] = 'Angeleno'
$data[14,11] = 'Segunda'
$data[14,12] = 16
$data[14,13] = 150000
$data[14,14] = 160000
$data[14,15] = 155000
$data[14,16] = '$/bins (450 kilos)'
$data[14,17] = 'Región Metropolitana'
$data[14,18] = 344
$data[14,19] = 450
$data[15,0] = 2
$data[15,1] = 'Comercializadora del Agro de Limarí'
$data[15,2] = 'Coquimbo'
$data[15,3] = 44721
$data[15,4] = 4
$data[15,5] = 'Fruta'
$data[15,6] = 100103
$data[15,7] = 'Frutos de hueso (carozo)'
$data[15,8] = 100103002
$data[15,9] = 'Ciruela'
$data[15,10] = 'Angeleno'
$data[15,11] = 'Primera'
$data[15,12] = 16
$data[15,13] = 190000
$data[15,14] = 200000
$data[15,15] = 195000
$data[15,16] = '$/bins (450 kilos)'
$data[15,17] = 'Región de O''Higgins'
$data[15,18] = 433
$data[15,19] = 450
$data[16,0] = 2
$data[16,1] = 'Comercializadora del Agro de Limarí'
$data[16,2] = 'Coquimbo'
$data[16,3] = 44721
$data[16,4] = 4
$data[16,5] = 'Fruta'
$data[16,6] = 100103
$data[16,7] = 'Frutos de hueso (carozo)'
$data[16,8] = 100103002
$data[16,9] = 'Ciruela'
$data[16,10] = 'Angeleno'
$data[16,11] = 'Segunda'
$data[16,12] = 16
$data[16,13] = 150000
$data[16,14] = 160000
$data[16,15] = 155000
$data[16,16] = '$/bins (450 kilos)'
$data[16,17] = 'Región de O''Higgins'
$data[16,18] = 344
$data[16,19] = 450
$data[17,0] = 2
$data[17,1] = 'Comercializadora del Agro de Limarí'
$data[17,2] = 'Coquimbo'
$data[17,3] = 44714
$data[17,4] = 4
$data[17,5] = 'Fruta'
$data[17,6] = 100103
$data[17,7] = 'Frutos de hueso (carozo)'
$data[17,8] = 100103002
$data[17,9] = 'Ciruela'
$data[17,10] = 'Angeleno'
$data[17,11] = 'Especial'
$data[17,12] = 10
$data[17,13] = 220000
$data[17,14] = 230000
$data[17,15] = 225000
$data[17,16] = '$/bins (450 kilos)'
$data[17,17] = 'Región de O''Higgins'
$data[17,18] = 500
$data[17,19] = 450
$data[18,0] = 2
$data[18,1] = 'Comercializadora del Agro de Limarí'
$data[18,2] = 'Coquimbo'
$data[18,3] = 44714
$data[18,4] = 4
$data[18,5] = 'Fruta'
$data[18,6] = 100103
$data[18,7] = 'Frutos de hueso (carozo)'
$data[18,8] = 100103002
$data[18,9] = 'Ciruela'
$data[18,10] = 'Angeleno'
$data[18,11] = 'Primera'
$data[18,12] = 16
$data[18,13] = 200000
$data[18,14] = 210000
$data[18,15] = 205000
$data[18,16] = '$/bins (450 kilos)'
$data[18,17] = 'Región de O''Higgins'
$data[18,18] = 456
$data[18,19] = 450
$data[19,0] = 2
$data[19,1] = 'Comercializadora del Agro de Limarí'
$data[19,2] = 'Coquimbo'
$data[19,3] = 44714
$data[19,4] = 4
$data[19,5] = 'Fruta'
$data[19,6] = 100103
$data[19,7] = 'Frutos de hueso (carozo)'
$data[19,8] = 100103002
$data[19,9] = 'Ciruela'
$data[19,10] = 'Angeleno'
$data[19,11] = 'Segunda'
$data[19,12] = 16
$data[19,13] = 170000
$data[19,14] = 180000
$data[19,15] = 175000
$data[19,16] = '$/bins (450 kilos)'
$data[19,17] = 'Región de O''Higgins'
$data[19,18] = 389
$data[19,19] = 450
$data[20,0] = 2
$data[20,1] = 'Comercializadora del Agro de Limarí'
$data[20,2] = 'Coquimbo'
$data[20,3] = 44672
$data[20,4] = 4
$data[20,5] = 'Fruta'
$data[20,6] = 100103
$data[20,7] = 'Frutos de hueso (carozo)'
$data[20,8] = 100103002
$data[20,9] = 'Ciruela'
$data[20,10] = 'Angeleno'
$data[20,11] = 'Primera'
$data[20,12] = 16
$data[20,13] = 190000
$data[20,14] = 200000
$data[20,15] = 195000
$data[20,16] = '$/bins (450 kilos)'
$data[20,17] = 'Región de O''Higgins'
$data[20,18] = 433
$data[20,19] = 450
$data[21,0] = 2
$data[21,1] = 'Comercializadora del Agro de Limarí'
$data[21,2] = 'Coquimbo'
$data[21,3] = 44672
$data[21,4] = 4
$data[21,5] = 'Fruta'
$data[21,6] = 100103
$data[21,7] = 'Frutos de hueso (carozo)'
$data[21,8] = 100103002
$data[21,9] = 'Ciruela'
$data[21,10] = 'Angeleno'
$data[21,11] = 'Segunda'
$data[21,12] = 16
$data[21,13] = 150000
$data[21,14] = 160000
$data[21,15] = 155000
$data[21,16] = '$/bins (450 kilos)'
$data[21,17] = 'Región de O''Higgins'
$data[21,18] = 344
$data[21,19] = 450
$data[22,0] = 2
$data[22,1] = 'Comercializadora del Agro de Limarí'
$data[22,2] = 'Coquimbo'
$data[22,3] = 44692
$data[22,4] = 4
$data[22,5] = 'Fruta'
$data[22,6] = 100103
$data[22,7] = 'Frutos de hueso (carozo)'
$data[22,8] = 100103002
$data[22,9] = 'Ciruela'
$data[22,10] = 'Angeleno'
$data[22,11] = 'Especial'
$data[22,12] = 20
$data[22,13] = 220000
$data[22,14] = 230000
$data[22,15] = 225000
$data[22,16] = '$/bins (450 kilos)'
$data[22,17] = 'Región de O''Higgins'
$data[22,18] = 500
$data[22,19] = 450
$data[23,0] = 2
$data[23,1] = 'Comercializadora del Agro de Limarí'
$data[23,2] = 'Coquimbo'
$data[23,3] = 44692
$data[23,4] = 4
$data[23,5] = 'Fruta'
$data[23,6] = 100103
$data[23,7] = 'Frutos de hueso (carozo)'
$data[23,8] = 100103002
$data[23,9] = 'Ciruela'
$data[23,10] = 'Angeleno'
$data[23,11] = 'Primera'
$data[23,12] = 20
$data[23,13] = 190000
$data[23,14] = 200000
$data[23,15] = 195000
$data[23,16] = '$/bins (450 kilos)'
$data[23,17] = 'Región de O''Higgins'
$data[23,18] = 433
$data[23,19] = 450
$data[24,0] = 2
$data[24,1] = 'Comercializadora del Agro de Limarí'
$data[24,2] = 'Coquimbo'
$data[24,3] = 44692
$data[24,4] = 4
$data[24,5] = 'Fruta'
$data[24,6] = 100103
$data[24,7] = 'Frutos de hueso (carozo)'
$data[24,8] = 100103002
$data[24,9] = 'Ciruela'
$data[24,10] = 'Angeleno'
$data[24,11] = 'Segunda'
$data[24,12] = 16
$data[24,13] = 150000
$data[24,14] = 160000
$data[24,15] = 155000
$data[24,16] = '$/bins (450 kilos)'
$data[24,17] = 'Región de O''Higgins'
$data[24,18] = 344
$data[24,19] = 450
$data[25,0] = 2
$data[25,1] = 'Comercializadora del Agro de Limarí'
$data[25,2] = 'Coquimbo'
$data[25,3] = 44629
$data[25,4] = 4
$data[25,5] = 'Fruta'
$data[25,6] = 100103
$data[25,7] = 'Frutos de hueso (carozo)'
$data[25,8] = 100103002
$data[25,9] = 'Ciruela'
$data[25,10] = 'Black Amber'
$data[25,11] = 'Primera'
$data[25,12] = 20
$data[25,13] = 220000
$data[25,14] = 230000
$data[25,15] = 225000
$data[25,16] = '$/bins (450 kilos)'
$data[25,17] = 'Región de O''Higgins'
$data[25,18] = 500
$data[25,19] = 450
$data[26,0] = 2
$data[26,1] = 'Comercializadora del Agro de Limarí'
$data[26,2] = 'Coquimbo'
$data[26,3] = 44629
$data[26,4] = 4
$data[26,5] = 'Fruta'
$data[26,6] = 100103
$data[26,7] = 'Frutos de hueso (carozo)'
$data[26,8] = 100103002
$data[26,9] = 'Ciruela'
$data[26,10] = 'Larry Ann'
$data[26,11] = 'Primera'
$data[26,12] = 20
$data[26,13] = 220000
$data[26,14] = 230000
$data[26,15] = 225000
$data[26,16] = '$/bins (450 kilos)'
$data[26,17] = 'Región de O''Higgins'
$data[26,18] = 500
$data[26,19] = 450
$data[27,0] = 2
$data[27,1] = 'Comercializadora del Agro de Limarí'
$data[27,2] = 'Coquimbo'
$data[27,3] = 44602
$data[27,4] = 4
$data[27,5] = 'Fruta'
$data[27,6] = 100103
$data[27,7] = 'Frutos de hueso (carozo)'
$data[27,8] = 100103002
$data[27,9] = 'Ciruela'
$data[27,10] = 'Larry Ann'
$data[27,11] = 'Primera'
$data[27,12] = 10
$data[27,13] = 185000
$data[27,14] = 190000
$data[27,15] = 187500
$data[27,16] = '$/bins (450 kilos)'
$data[27,17] = 'Región Metropolitana'
$data[27,18] = 417
$data[27,19] = 450
$data[28,0] = 2
$data[28,1] = 'Comercializadora del Agro de Limarí'
$data[28,2] = 'Coquimbo'
$data[28,3] = 44602
$data[28,4] = 4
$data[28,5] = 'Fruta'
$data[28,6] = 100103
$data[28,7] = 'Frutos de hueso (carozo)'
$data[28,8] = 100103002
$data[28,9] = 'Ciruela'
$data[28,10] = 'Larry Ann'
$data[28,11] = 'Segunda'
$data[28,12] = 10
$data[28,13] = 155000
$data[28,14] = 160000
$data[28,15] = 157500
$data[28,16] = '$/bins (450 kilos)'
$data[28,17] = 'Región Metropolitana'
$data[28,18] = 350
$data[28,19] = 450
$data[29,0] = 2
$data[29,1] = 'Comercializadora del Agro de Limarí'
$data[29,2] = 'Coquimbo'
$data[29,3] = 44707
$data[29,4] = 4
$data[29,5] = 'Fruta'
$data[29,6] = 100103
$data[29,7] = 'Frutos de hueso (carozo)'
$data[29,8] = 100103002
$data[29,9] = 'Ciruela'
$data[29,10] = 'Angeleno'
$data[29,11] = 'Primera'
$data[29,12] = 16
$data[29,13] = 190000
$data[29,14] = 200000
$data[29,15] = 195000
$data[29,16] = '$/bins (450 kilos)'
$data[29,17] = 'Región de O''Higgins'
$data[29,18] = 433
$data[29,19] = 450
$data[30,0] = 2
$data[30,1] = 'Comercializadora del Agro de Limarí'
$data[30,2] = 'Coquimbo'
$data[30,3] = 44707
$data[30,4] = 4
$data[30,5] = 'Fruta'
$data[30,6] = 100103
$data[30,7] = 'Frutos de hueso (carozo)'
$data[30,8] = 100103002
$data[30,9] = 'Ciruela'
$data[30,10] = 'Angeleno'
$data[30,11] = 'Segunda'
$data[30,12] = 20
$data[30,13] = 150000
$data[30,14] = 160000
$data[30,15] = 155000
$data[30,16] = '$/bins (450 kilos)'
$data[30,17] = 'Región de O''Higgins'
$data[30,18] = 344
$data[30,19] = 450
$data[31,0] = 2
$data[31,1] = 'Comercializadora del Agro de Limarí'
$data[31,2] = 'Coquimbo'
$data[31,3] = 44209
$data[31,4] = 4
$data[31,5] = 'Fruta'
$data[31,6] = 100103
$data[31,7] = 'Frutos de hueso (carozo)'
$data[31,8] = 100103002
$data[31,9] = 'Ciruela'
$data[31,10] = 'Black Amber'
$data[31,11] = 'Primera'
$data[31,12] = 300
$data[31,13] = 15500
$data[31,14] = 16000
$data[31,15] = 15750
$data[31,16] = '$/caja 16 kilos granel'
$data[31,17] = 'Región Metropolitana'
$data[31,18] = 984
$data[31,19] = 16
$data[32,0] = 2
$data[32,1] = 'Comercializadora del Agro de Limarí'
$data[32,2] = 'Coquimbo'
$data[32,3] = 44951
$data[32,4] = 4
$data[32,5] = 'Fruta'
$data[32,6] = 100103
$data[32,7] = 'Frutos de hueso (carozo)'
$data[32,8] = 100103002
$data[32,9] = 'Ciruela'
$data[32,10] = 'Black Amber'
$data[32,11] = 'Primera'
$data[32,12] = 16
$data[32,13] = 300000
$data[32,14] = 310000
$data[32,15] = 305000
$data[32,16] = '$/bins (450 kilos)'
$data[32,17] = 'Región de O''Higgins'
$data[32,18] = 678
$data[32,19] = 450
$data[33,0] = 2
$data[33,1] = 'Comercializadora del Agro de Limarí'
$data[33,2] = 'Coquimbo'
$data[33,3] = 44951
$data[33,4] = 4
$data[33,5] = 'Fruta'
$data[33,6] = 100103
$data[33,7] = 'Frutos de hueso (carozo)'
$data[33,8] = 100103002
$data[33,9] = 'Ciruela'
$data[33,10] = 'Black Amber'
$data[33,11] = 'Segunda'
$data[33,12] = 20
$data[33,13] = 230000
$data[33,14] = 240000
$data[33,15] = 235000
$data[33,16] = '$/bins (450 kilos)'
$data[33,17] = 'Región de O''Higgins'
$data[33,18] = 522
$data[33,19] = 450
$data[34,0] = 2
$data[34,1] = 'Comercializadora del Agro de Limarí'
$data[34,2] = 'Coquimbo'
$data[34,3] = 44658
$data[34,4] = 4
$data[34,5] = 'Fruta'
$data[34,6] = 100103
$data[34,7] = 'Frutos de hueso (carozo)'
$data[34,8] = 100103002
$data[34,9] = 'Ciruela'
$data[34,10] = 'Angeleno'
$data[34,11] = 'Especial'
$data[34,12] = 16
$data[34,13] = 220000
$data[34,14] = 230000
$data[34,15] = 225000
$data[34,16] = '$/bins (450 kilos)'
$data[34,17] = 'Región de O''Higgins'
$data[34,18] = 500
$data[34,19] = 450
$data[35,0] = 2
$data[35,1] = 'Comercializadora del Agro de Limarí'
$data[35,2] = 'Coquimbo'
$data[35,3] = 44658
$data[35,4] = 4
$data[35,5] = 'Fruta'
$data[35,6] = 100103
$data[35,7] = 'Frutos de hueso (carozo)'
$data[35,8] = 100103002
$data[35,9] = 'Ciruela'
$data[35,10] = 'Angeleno'
$data[35,11] = 'Primera'
$data[35,12] = 16
$data[35,13] = 180000
$data[35,14] = 190000
$data[35,15] = 185000
$data[35,16] = '$/bins (450 kilos)'
$data[35,17] = 'Región de O''Higgins'
$data[35,18] = 411
$data[35,19] = 450
$data[36,0] = 2
$data[36,1] = 'Comercializadora del Agro de Limarí'
$data[36,2] = 'Coquimbo'
$data[36,3] = 44658
$data[36,4] = 4
$data[36,5] = 'Fruta'
$data[36,6] = 100103
$data[36,7] = 'Frutos de hueso (carozo)'
$data[36,8] = 100103002
$data[36,9] = 'Ciruela'
$data[36,10] = 'Angeleno'
$data[36,11] = 'Segunda'
$data[36,12] = 16
$data[36,13] = 150000
$data[36,14] = 160000
$data[36,15] = 155000
$data[36,16] = '$/bins (450 kilos)'
$data[36,17] = 'Región de O''Higgins'
$data[36,18] = 344
$data[36,19] = 450
$data[37,0] = 2
$data[37,1] = 'Comercializadora del Agro de Limarí'
$data[37,2] = 'Coquimbo'
$data[37,3] = 44644
$data[37,4] = 4
$data[37,5] = 'Fruta'
$data[37,6] = 100103
$data[37,7] = 'Frutos de hueso (carozo)'
$data[37,8] = 100103002
$data[37,9] = 'Ciruela'
$data[37,10] = 'Angeleno'
$data[37,11] = 'Especial'
$data[37,12] = 10
$data[37,13] = 230000
$data[37,14] = 240000
$data[37,15] = 235000
$data[37,16] = '$/bins (450 kilos)'
$data[37,17] = 'Región de O''Higgins'
$data[37,18] = 522
$data[37,19] = 450
$data[38,0] = 2
$data[38,1] = 'Comercializadora del Agro de Limarí'
$data[38,2] = 'Coquimbo'
$data[38,3] = 44644
$data[38,4] = 4
$data[38,5] = 'Fruta'
$data[38,6] = 100103
$data[38,7] = 'Frutos de hueso (carozo)'
$data[38,8] = 100103002
$data[38,9] = 'Ciruela'
$data[38,10] = 'Angeleno'
$data[38,11] = 'Primera'
$data[38,12] = 16
$data[38,13] = 210000
$data[38,14] = 220000
$data[38,15] = 215000
$data[38,16] = '$/bins (450 kilos)'
$data[38,17] = 'Región de O''Higgins'
$data[38,18] = 478
$data[38,19] = 450
$data[39,0] = 2
$data[39,1] = 'Comercializadora del Agro de Limarí'
$data[39,2] = 'Coquimbo'
$data[39,3] = 44644
$data[39,4] = 4
$data[39,5] = 'Fruta'
$data[39,6] = 100103
$data[39,7] = 'Frutos de hueso (carozo)'
$data[39,8] = 100103002
$data[39,9] = 'Ciruela'
$data[39,10] = 'Angeleno'
$data[39,11] = 'Segunda'
$data[39,12] = 16
$data[39,13] = 150000
$data[39,14] = 160000
$data[39,15] = 155000
$data[39,16] = '$/bins (450 kilos)'
$data[39,17] = 'Región de O''Higgins'
$data[39,18] = 344
$data[39,19] = 450
$data[40,0] = 2
$data[40,1] = 'Comercializadora del Agro de Limarí'
$data[40,2] = 'Coquimbo'
$data[40,3] = 44224
$data[40,4] = 4
$data[40,5] = 'Fruta'
$data[40,6] = 100103
$data[40,7] = 'Frutos de hueso (carozo)'
$data[40,8] = 100103002
$data[40,9] = 'Ciruela'
$data[40,10] = 'Black Amber'
$data[40,11] = 'Especial'
$data[40,12] = 100
$data[40,13] = 16500
$data[40,14] = 17000
$data[40,15] = 16750
$data[40,16] = '$/caja 16 kilos granel'
$data[40,17] = 'Región Metropolitana'
$data[40,18] = 1047
$data[40,19] = 16
$data[41,0] = 2
$data[41,1] = 'Comercializadora del Agro de Limarí'
$data[41,2] = 'Coquimbo'
$data[41,3] = 44224
$data[41,4] = 4
$data[41,5] = 'Fruta'
$data[41,6] = 100103
$data[41,7] = 'Frutos de hueso (carozo)'
$data[41,8] = 100103002
$data[41,9] = 'Ciruela'
$data[41,10] = 'Black Amber'
$data[41,11] = 'Primera'
$data[41,12] = 200
$data[41,13] = 14500
$data[41,14] = 15000
$data[41,15] = 14750
$data[41,16] = '$/caja 16 kilos granel'
$data[41,17] = 'Región Metropolitana'
$data[41,18] = 922
$data[41,19] = 16
$data[42,0] = 2
$data[42,1] = 'Comercializadora del Agro de Limarí'
$data[42,2] = 'Coquimbo'
$data[42,3] = 44224
$data[42,4] = 4
$data[42,5] = 'Fruta'
$data[42,6] = 100103
$data[42,7] = 'Frutos de hueso (carozo)'
$data[42,8] = 100103002
$data[42,9] = 'Ciruela'
$data[42,10] = 'Black Amber'
$data[42,11] = 'Segunda'
$data[42,12] = 200
$data[42,13] = 12500
$data[42,14] = 13000
$data[42,15] = 12750
$data[42,16] = '$/caja 16 kilos granel'
$data[42,17] = 'Región Metropolitana'
$data[42,18] = 797
$data[42,19] = 16
$data[43,0] = 2
$data[43,1] = 'Comercializadora del Agro de Limarí'
$data[43,2] = 'Coquimbo'
$data[43,3] = 44637
$data[43,4] = 4
$data[43,5] = 'Fruta'
$data[43,6] = 100103
$data[43,7] = 'Frutos de hueso (carozo)'
$data[43,8] = 100103002
$data[43,9] = 'Ciruela'
$data[43,10] = 'Angeleno'
$data[43,11] = 'Especial'
$data[43,12] = 20
$data[43,13] = 255000
$data[43,14] = 260000
$data[43,15] = 257500
$data[43,16] = '$/bins (450 kilos)'
$data[43,17] = 'Región Metropolitana'
$data[43,18] = 572
$data[43,19] = 450
$data[44,0] = 2
$data[44,1] = 'Comercializadora del Agro de Limarí'
$data[44,2] = 'Coquimbo'
$data[44,3] = 44637
$data[44,4] = 4
$data[44,5] = 'Fruta'
$data[44,6] = 100103
$data[44,7] = 'Frutos de hueso (carozo)'
$data[44,8] = 100103002
$data[44,9] = 'Ciruela'
$data[44,10] = 'Angeleno'
$data[44,11] = 'Primera'
$data[44,12] = 20
$data[44,13] = 225000
$data[44,14] = 230000
$data[44,15] = 227500
$data[44,16] = '$/bins (450 kilos)'
$data[44,17] = 'Región Metropolitana'
$data[44,18] = 506
$data[44,19] = 450
$data[45,0] = 2
$data[45,1] = 'Comercializadora del Agro de Limarí'
$data[45,2] = 'Coquimbo'
$data[45,3] = 44616
$data[45,4] = 4
$data[45,5] = 'Fruta'
$data[45,6] = 100103
$data[45,7] = 'Frutos de hueso (carozo)'
$data[45,8] = 100103002
$data[45,9] = 'Ciruela'
$data[45,10] = 'Angeleno'
$data[45,11] = 'Especial'
$data[45,12] = 10
$data[45,13] = 220000
$data[45,14] = 230000
$data[45,15] = 225000
$data[45,16] = '$/bins (450 kilos)'
$data[45,17] = 'Región de O''Higgins'
$data[45,18] = 500
$data[45,19] = 450
$data[46,0] = 2
$data[46,1] = 'Comercializadora del Agro de Limarí'
$data[46,2] = 'Coquimbo'
$data[46,3] = 44616
$data[46,4] = 4
$data[46,5] = 'Fruta'
$data[46,6] = 100103
$data[46,7] = 'Frutos de hueso (carozo)'
$data[46,8] = 100103002
$data[46,9] = 'Ciruela'
$data[46,10] = 'Angeleno'
$data[46,11] = 'Primera'
$data[46,12] = 16
$data[46,13] = 200000
$data[46,14] = 210000
$data[46,15] = 205000
$data[46,16] = '$/bins (450 kilos)'
$data[46,17] = 'Región de O''Higgins'
$data[46,18] = 456
$data[46,19] = 450
$data[47,0] = 2
$data[47,1] = 'Comercializadora del Agro de Limarí'
$data[47,2] = 'Coquimbo'
$data[47,3] = 44616
$data[47,4] = 4
$data[47,5] = 'Fruta'
$data[47,6] = 100103
$data[47,7] = 'Frutos de hueso (carozo)'
$data[47,8] = 100103002
$data[47,9] = 'Ciruela'
$data[47,10] = 'Angeleno'
$data[47,11] = 'Segunda'
$data[47,12] = 10
$data[47,13] = 160000
$data[47,14] = 170000
$data[47,15] = 165000
$data[47,16] = '$/bins (450 kilos)'
$data[47,17] = 'Región de O''Higgins'
$data[47,18] = 367
$data[47,19] = 450
$data[48,0] = 2
$data[48,1] = 'Comercializadora del Agro de Limarí'
$data[48,2] = 'Coquimbo'
$data[48,3] = 44616
$data[48,4] = 4
$data[48,5] = 'Fruta'
$data[48,6] = 100103
$data[48,7] = 'Frutos de hueso (carozo)'
$data[48,8] = 100103002
$data[48,9] = 'Ciruela'
$data[48,10] = 'Black Amber'
$data[48,11] = 'Primera'
$data[48,12] = 20
$data[48,13] = 200000
$data[48,14] = 210000
$data[48,15] = 205000
$data[48,16] = '$/bins (450 kilos)'
$data[48,17] = 'Región Metropolitana'
$data[48,18] = 456
$data[48,19] = 450
$data[49,0] = 2
$data[49,1] = 'Comercializadora del Agro de Limarí'
$data[49,2] = 'Coquimbo'
$data[49,3] = 44616
$data[49,4] = 4
$data[49,5] = 'Fruta'
$data[49,6] = 100103
$data[49,7] = 'Frutos de hueso (carozo)'
$data[49,8] = 100103002
$data[49,9] = 'Ciruela'
$data[49,10] = 'Black Amber'
$data[49,11] = 'Segunda'
$data[49,12] = 20
$data[49,13] = 150000
$data[49,14] = 160000
$data[49,15] = 155000
$data[49,16] = '$/bins (450 kilos)'
$data[49,17] = 'Región Metropolitana'
$data[49,18] = 344
$data[49,19] = 450
$data[50,0] = 2
$data[50,1] = 'Comercializadora del Agro de Limarí'
$data[50,2] = 'Coquimbo'
$data[50,3] = 44643
$data[50,4] = 4
$data[50,5] = 'Fruta'
$data[50,6] = 100103
$data[50,7] = 'Frutos de hueso (carozo)'
$data[50,8] = 100103002
$data[50,9] = 'Ciruela'
$data[50,10] = 'Angeleno'
$data[50,11] = 'Especial'
$data[50,12] = 10
$data[50,13] = 230000
$data[50,14] = 240000
$data[50,15] = 235000
$data[50,16] = '$/bins (450 kilos)'
$data[50,17] = 'Región de O''Higgins'
$data[50,18] = 522
$data[50,19] = 450
$data[51,0] = 2
$data[51,1] = 'Comercializadora del Agro de Limarí'
$data[51,2] = 'Coquimbo'
$data[51,3] = 44643
$data[51,4] = 4
$data[51,5] = 'Fruta'
$data[51,6] = 100103
$data[51,7] = 'Frutos de hueso (carozo)'
$data[51,8] = 100103002
$data[51,9] = 'Ciruela'
$data[51,10] = 'Angeleno'
$data[51,11] = 'Primera'
$data[51,12] = 16
$data[51,13] = 210000
$data[51,14] = 220000
$data[51,15] = 215000
$data[51,16] = '$/bins (450 kilos)'
$data[51,17] = 'Región de O''Higgins'
$data[51,18] = 478
$data[51,19] = 450
$data[52,0] = 2
$data[52,1] = 'Comercializadora del Agro de Limarí'
$data[52,2] = 'Coquimbo'
$data[52,3] = 44671
$data[52,4] = 4
$data[52,5] = 'Fruta'
$data[52,6] = 100103
$data[52,7] = 'Frutos de hueso (carozo)'
$data[52,8] = 100103002
$data[52,9] = 'Ciruela'
$data[52,10] = 'Angeleno'
$data[52,11] = 'Primera'
$data[52,12] = 16
$data[52,13] = 190000
$data[52,14] = 200000
$data[52,15] = 195000
$data[52,16] = '$/bins (450 kilos)'
$data[52,17] = 'Región de O''Higgins'
$data[52,18] = 433
$data[52,19] = 450
$data[53,0] = 2
$data[53,1] = 'Comercializadora del Agro de Limarí'
$data[53,2] = 'Coquimbo'
$data[53,3] = 44671
$data[53,4] = 4
$data[53,5] = 'Fruta'
$data[53,6] = 100103
$data[53,7] = 'Frutos de hueso (carozo)'
$data[53,8] = 100103002
$data[53,9] = 'Ciruela'
$data[53,10] = 'Angeleno'
$data[53,11] = 'Segunda'
$data[53,12] = 10
$data[53,13] = 150000
$data[53,14] = 160000
$data[53,15] = 155000
$data[53,16] = '$/bins (450 kilos)'
$data[53,17] = 'Región de O''Higgins'
$data[53,18] = 344
$data[53,19] = 450
$data[54,0] = 2
$data[54,1] = 'Comercializadora del Agro de Limarí'
$data[54,2] = 'Coquimbo'
$data[54,3] = 44657
$data[54,4] = 4
$data[54,5] = 'Fruta'
$data[54,6] = 100103
$data[54,7] = 'Frutos de hueso (carozo)'
$data[54,8] = 100103002
$data[54,9] = 'Ciruela'
$data[54,10] = 'Angeleno'
$data[54,11] = 'Primera'
$data[54,12] = 16
$data[54,13] = 180000
$data[54,14] = 190000
$data[54,15] = 185000
$data[54,16] = '$/bins (450 kilos)'
$data[54,17] = 'Región de O''Higgins'
$data[54,18] = 411
$data[54,19] = 450
$data[55,0] = 2
$data[55,1] = 'Comercializadora del Agro de Limarí'
$data[55,2] = 'Coquimbo'
$data[55,3] = 44636
$data[55,4] = 4
$data[55,5] = 'Fruta'
$data[55,6] = 100103
$data[55,7] = 'Frutos de hueso (carozo)'
$data[55,8] = 100103002
$data[55,9] = 'Ciruela'
$data[55,10] = 'Angeleno'
$data[55,11] = 'Primera'
$data[55,12] = 16
$data[55,13] = 235000
$data[55,14] = 240000
$data[55,15] = 237500
$data[55,16] = '$/bins (450 kilos)'
$data[55,17] = 'Región de O''Higgins'
$data[55,18] = 528
$data[55,19] = 450
$data[56,0] = 2
$data[56,1] = 'Comercializadora del Agro de Limarí'
$data[56,2] = 'Coquimbo'
$data[56,3] = 44636
$data[56,4] = 4
$data[56,5] = 'Fruta'
$data[56,6] = 100103
$data[56,7] = 'Frutos de hueso (carozo)'
$data[56,8] = 100103002
$data[56,9] = 'Ciruela'
$data[56,10] = 'Angeleno'
$data[56,11] = 'Segunda'
$data[56,12] = 20
$data[56,13] = 185000
$data[56,14] = 190000
$data[56,15] = 187500
$data[56,16] = '$/bins (450 kilos)'
$data[56,17] = 'Región de O''Higgins'
$data[56,18] = 417
$data[56,19] = 450
$data[57,0] = 2
$data[57,1] = 'Comercializadora del Agro de Limarí'
$data[57,2] = 'Coquimbo'
$data[57,3] = 44595
$data[57,4] = 4
$data[57,5] = 'Fruta'
$data[57,6] = 100103
$data[57,7] = 'Frutos de hueso (carozo)'
$data[57,8] = 100103002
$data[57,9] = 'Ciruela'
$data[57,10] = 'Black Amber'
$data[57,11] = 'Especial'
$data[57,12] = 160
$data[57,13] = 15500
$data[57,14] = 16000
$data[57,15] = 15750
$data[57,16] = '$/caja 15 kilos granel'
$data[57,17] = 'Región de O''Higgins'
$data[57,18] = 1050
$data[57,19] = 15
$data[58,0] = 2
$data[58,1] = 'Comercializadora del Agro de Limarí'
$data[58,2] = 'Coquimbo'
$data[58,3] = 44595
$data[58,4] = 4
$data[58,5] = 'Fruta'
$data[58,6] = 100103
$data[58,7] = 'Frutos de hueso (carozo)'
$data[58,8] = 100103002
$data[58,9] = 'Ciruela'
$data[58,10] = 'Black Amber'
$data[58,11] = 'Primera'
$data[58,12] = 200
$data[58,13] = 13500
$data[58,14] = 14000
$data[58,15] = 13750
$data[58,16] = '$/caja 15 kilos granel'
$data[58,17] = 'Región de O''Higgins'
$data[58,18] = 917
$data[58,19] = 15
$data[59,0] = 2
$data[59,1] = 'Comercializadora del Agro de Limarí'
$data[59,2] = 'Coquimbo'
$data[59,3] = 44595
$data[59,4] = 4
$data[59,5] = 'Fruta'
$data[59,6] = 100103
$data[59,7] = 'Frutos de hueso (carozo)'
$data[59,8] = 100103002
$data[59,9] = 'Ciruela'
$data[59,10] = 'Black Amber'
$data[59,11] = 'Segunda'
$data[59,12] = 300
$data[59,13] = 10500
$data[59,14] = 11000
$data[59,15] = 10750
$data[59,16] = '$/caja 15 kilos granel'
$data[59,17] = 'Región de O''Higgins'
$data[59,18] = 717
$data[59,19] = 15

$startRow = 27
$endRow = 86
$rng = $ws.Range($ws.Cells.Item($startRow,1), $ws.Cells.Item($endRow,20))
$rng.Value = $data

# Ensure the date column (D) uses the correct date number format for
# the newly appended rows (85 and 86), which previously had no style.
$ws.Range('D85:D86').NumberFormat = 'YYYY-MM-DD HH:MM:SS'

